# "Add files via upload" — refresh the private-display test-case log:
# swap test case #1's "Given Input" sequence from the old
# operator-overloading / copy-constructor / member-function flow to the
# new default-constructor / assignment-operator-overloading flow (and the
# matching re-run of that same flow for test case #3), widen columns D/E
# to fit the new (longer) text, and move the visible window down to the
# newly-added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Test case #1 (rows 2-12): new "Given Input" step sequence ---------
$ws.Range("D2").Value  = "using default constructor"
$ws.Range("E2").Value  = "using default constructor"
$ws.Range("D3").Value  = "no=0"
$ws.Range("E3").Value  = "no=0"
$ws.Range("D4").Value  = "name=--"
$ws.Range("E4").Value  = "name=--"
$ws.Range("D5").Value  = "using parameterized constructor"
$ws.Range("E5").Value  = "using parameterized constructor"
$ws.Range("D6").Value  = "no=10"
$ws.Range("E6").Value  = "no=10"
$ws.Range("D7").Value  = "name=prasuna"
$ws.Range("E7").Value  = "name=prasuna"
$ws.Range("D8").Value  = "After copy contructor"
$ws.Range("E8").Value  = "After copy contructor"
$ws.Range("D9").Value  = "no=10 name=prasuna"
$ws.Range("E9").Value  = "no=10 name=prasuna"
$ws.Range("D10").Value = "no=9 name=raj"
$ws.Range("E10").Value = "no=9 name=raj"
$ws.Range("D11").Value = "after assignment operator overloading"
$ws.Range("E11").Value = "after assignment operator overloading"
$ws.Range("D12").Value = "no=0 name=--"
$ws.Range("E12").Value = "no=0 name=--"

# --- Test case #3 (rows 17-27): same replay of the new sequence --------
$ws.Range("D17").Value = "using default constructor"
$ws.Range("E17").Value = "using default constructor"
$ws.Range("D18").Value = "no=0"
$ws.Range("E18").Value = "no=0"
$ws.Range("D19").Value = "name=--"
$ws.Range("E19").Value = "name=--"
$ws.Range("D20").Value = "using parameterized constructor"
$ws.Range("E20").Value = "using parameterized constructor"
$ws.Range("D21").Value = "no=10"
$ws.Range("E21").Value = "no=10"
$ws.Range("D22").Value = "name=prasuna"
$ws.Range("E22").Value = "name=prasuna"
$ws.Range("D23").Value = "After copy contructor"
$ws.Range("E23").Value = "After copy contructor"
$ws.Range("D24").Value = "no=10 name=prasuna"
$ws.Range("E24").Value = "no=10 name=prasuna"
$ws.Range("D25").Value = "no=0 name=--"
$ws.Range("E25").Value = "no=0 name=--"
$ws.Range("D26").Value = "after assignment operator overloading"
$ws.Range("E26").Value = "after assignment operator overloading"
$ws.Range("D27").Value = "no=10 name=prasuna"
$ws.Range("E27").Value = "no=10 name=prasuna"

# --- Widen D/E to fit the new, longer text ------------------------------
$ws.Columns.Item(4).ColumnWidth = 35.4518229166667
$ws.Columns.Item(5).ColumnWidth = 34.0221354166667

# --- Scroll the view down to the newly-added rows and reselect ---------
$ws.Range("E29").Select()
